# Hortaliza / Terminal Hortofrutícola Agro Chillán - Cebolla
# Insert 3 new weekly price rows right after the existing row for
# mercado id 7 / Ñuble (pushes everything from old row 559 onward down
# by 3 rows, growing the used range from A1:R639 to A1:R642).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the current row 559.
$ws.Rows("559:561").Insert()

# Full row data (columns A..R) for the 3 newly inserted rows.
$newData = @(
    @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44776, 16, 100112004, "Cebolla", "Sin especificar", "1a (guarda)", 120, 8000, 8500, 8250, "`$/malla 25 kilos", "Región del Maule", 330, 25, "Hortaliza"),
    @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44776, 16, 100112004, "Cebolla", "Sin especificar", "2a (guarda)", 120, 5000, 5000, 5000, "`$/malla 15 kilos", "Región del Maule", 333, 15, "Hortaliza"),
    @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44776, 16, 100112004, "Cebolla", "Sin especificar", "2a (guarda)", 200, 6500, 7000, 6750, "`$/malla 25 kilos", "Región del Maule", 270, 25, "Hortaliza")
)

$startRow = 559
for ($i = 0; $i -lt $newData.Length; $i++) {
    $rowValues = $newData[$i]
    for ($col = 0; $col -lt $rowValues.Length; $col++) {
        $ws.Cells.Item($startRow + $i, $col + 1).Value = $rowValues[$col]
    }
}
